$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1096.6666
$ws.Range("I96").Value = 1145
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 3435
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = -2062
$ws.Range("N96").Value = -5746

$ws.Range("H100").Value = 8674.5
$ws.Range("I100").Value = 804
$ws.Range("J100").Value = 10248.6
$ws.Range("K100").Value = 804
$ws.Range("L100").Value = 10248.6
$ws.Range("M100").Value = -263
$ws.Range("N100").Value = -11330.6

$ws.Range("H132").Value = 1666.04
$ws.Range("I132").Value = 1564.6136
$ws.Range("K132").Value = 4693.8408
$ws.Range("M132").Value = -2163.8408

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6187.6484
$ws.Range("I32").Value = 4739.6553
$ws.Range("K32").Value = 4739.6553
$ws.Range("M32").Value = -4452.6553

$ws.Range("H102").Value = 2457.625
$ws.Range("I102").Value = 2344.4285
$ws.Range("K102").Value = 2344.4285
$ws.Range("M102").Value = -722.4285

$ws.Range("H110").Value = 4033.2307
$ws.Range("I110").Value = 2627.818
$ws.Range("J110").Value = 11763
$ws.Range("K110").Value = 2627.818
$ws.Range("L110").Value = 11763
$ws.Range("M110").Value = -582.8180000000002
$ws.Range("N110").Value = -15853

$ws.Range("H122").Value = 3260.0435
$ws.Range("I122").Value = 2479.2
$ws.Range("J122").Value = 3860.6924
$ws.Range("K122").Value = 7437.599999999999
$ws.Range("L122").Value = 11582.0772
$ws.Range("M122").Value = -4987.599999999999
$ws.Range("N122").Value = -16482.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = ""

$ws.Range("H22").Value = 754.4286
$ws.Range("J22").Value = 1327.3334
$ws.Range("L22").Value = 1327.3334
$ws.Range("N22").Value = -1673.3334

$ws.Range("H94").Value = 2448
$ws.Range("I94").Value = 1688
$ws.Range("K94").Value = 1688
$ws.Range("M94").Value = -1237

$ws.Range("H105").Value = 21874.072
$ws.Range("I105").Value = 35623
$ws.Range("K105").Value = 35623
$ws.Range("M105").Value = -33876

$ws.Range("H110").Value = 37085.668
$ws.Range("J110").Value = 37085.668
$ws.Range("L110").Value = 37085.668
$ws.Range("N110").Value = -45265.668

$ws.Range("H130").Value = 78000
$ws.Range("J130").Value = 78000
$ws.Range("L130").Value = 78000
$ws.Range("N130").Value = -88040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29499.418
$ws.Range("I31").Value = 4072.3845
$ws.Range("K31").Value = 4072.3845
$ws.Range("M31").Value = -3777.3845

$ws.Range("H34").Value = 29499.418
$ws.Range("I34").Value = 4072.3845
$ws.Range("K34").Value = 4072.3845
$ws.Range("M34").Value = -3870.3845

$ws.Range("H58").Value = 4048.2856
$ws.Range("I58").Value = 1956.7
$ws.Range("J58").Value = 9277.25
$ws.Range("K58").Value = 1956.7
$ws.Range("L58").Value = 9277.25
$ws.Range("M58").Value = -1753.7
$ws.Range("N58").Value = -9683.25

$ws.Range("H86").Value = 8718.200000000001
$ws.Range("I86").Value = 8803
$ws.Range("K86").Value = 8803
$ws.Range("M86").Value = -7680

$ws.Range("H89").Value = 8718.200000000001
$ws.Range("I89").Value = 8803
$ws.Range("K89").Value = 44015
$ws.Range("M89").Value = -38399

$ws.Range("H105").Value = 1698.5588
$ws.Range("I105").Value = 929.8
$ws.Range("K105").Value = 929.8
$ws.Range("M105").Value = 817.2

$ws.Range("H132").Value = 2588.5518
$ws.Range("I132").Value = 1759.4166
$ws.Range("K132").Value = 5278.2498
$ws.Range("M132").Value = -2748.2498

$ws.Range("H135").Value = 68946.64
$ws.Range("J135").Value = 68946.64
$ws.Range("L135").Value = 68946.64
$ws.Range("N135").Value = -79086.64

$ws.Range("H136").Value = 4048.2856
$ws.Range("I136").Value = 1956.7
$ws.Range("J136").Value = 9277.25
$ws.Range("K136").Value = 5870.1
$ws.Range("L136").Value = 27831.75
$ws.Range("M136").Value = -3320.1
$ws.Range("N136").Value = -32931.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8882366
$ws.Range("I4").Value = 4000308.5
$ws.Range("K4").Value = 12000925.5
$ws.Range("M4").Value = -12000813.5

$ws.Range("H22").Value = 2639.5
$ws.Range("I22").Value = 979
$ws.Range("J22").Value = 4300
$ws.Range("K22").Value = 2937
$ws.Range("L22").Value = 12900
$ws.Range("M22").Value = -2768
$ws.Range("N22").Value = -13238

$ws.Range("H27").Value = 2639.5
$ws.Range("I27").Value = 979
$ws.Range("J27").Value = 4300
$ws.Range("K27").Value = 2937
$ws.Range("L27").Value = 12900
$ws.Range("M27").Value = -2835
$ws.Range("N27").Value = -13104

$ws.Range("H40").Value = 129.16667
$ws.Range("I40").Value = 62.5
$ws.Range("J40").Value = 162.5
$ws.Range("K40").Value = 250
$ws.Range("L40").Value = 650
$ws.Range("M40").Value = -181
$ws.Range("N40").Value = -788

$ws.Range("H42").Value = 14502
$ws.Range("J42").Value = 14502
$ws.Range("L42").Value = 43506
$ws.Range("N42").Value = -44574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 188332.5

$ws.Range("H30").Value = 188332.5

$ws.Range("H32").Value = 48000
$ws.Range("J32").Value = 48000
$ws.Range("L32").Value = 48000
$ws.Range("N32").Value = -48592

$ws.Range("H70").Value = 16736.309
$ws.Range("J70").Value = 19113.555
$ws.Range("L70").Value = 19113.555
$ws.Range("N70").Value = -19653.555

$ws.Range("H73").Value = 16736.309
$ws.Range("J73").Value = 19113.555
$ws.Range("L73").Value = 19113.555
$ws.Range("N73").Value = -20985.555

$ws.Range("H113").Value = 2607.5
$ws.Range("I113").Value = 1829.3529
$ws.Range("K113").Value = 1829.3529
$ws.Range("M113").Value = 340.6470999999999

$ws.Range("H132").Value = 5562.392
$ws.Range("I132").Value = 5293.64
$ws.Range("J132").Value = 19000
$ws.Range("K132").Value = 15880.92
$ws.Range("L132").Value = 57000
$ws.Range("M132").Value = -13350.92
$ws.Range("N132").Value = -62060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7115.8823
$ws.Range("I7").Value = 5766.5137
$ws.Range("K7").Value = 5766.5137
$ws.Range("M7").Value = -5654.5137

$ws.Range("H22").Value = 5220.2
$ws.Range("I22").Value = 2738.077
$ws.Range("K22").Value = 2738.077
$ws.Range("M22").Value = -2443.077

$ws.Range("H27").Value = 5220.2
$ws.Range("I27").Value = 2738.077
$ws.Range("K27").Value = 2738.077
$ws.Range("M27").Value = -2631.077

$ws.Range("H31").Value = 9182.833000000001
$ws.Range("J31").Value = 15550
$ws.Range("L31").Value = 15550
$ws.Range("N31").Value = -16046

$ws.Range("H46").Value = 3728.16
$ws.Range("I46").Value = 2543.2856
$ws.Range("J46").Value = 4188.9443
$ws.Range("K46").Value = 2543.2856
$ws.Range("L46").Value = 4188.9443
$ws.Range("M46").Value = -2355.2856
$ws.Range("N46").Value = -4564.9443

$ws.Range("H61").Value = 3622.5925
$ws.Range("I61").Value = 2560
$ws.Range("K61").Value = 2560
$ws.Range("M61").Value = -2358

$ws.Range("H113").Value = 3622.5925
$ws.Range("I113").Value = 2560
$ws.Range("K113").Value = 2560
$ws.Range("M113").Value = -390

$ws.Range("H126").Value = 7115.8823
$ws.Range("I126").Value = 5766.5137
$ws.Range("K126").Value = 17299.5411
$ws.Range("M126").Value = -14829.5411

$ws.Range("H132").Value = 6161.5435
$ws.Range("I132").Value = 5414.552
$ws.Range("J132").Value = 7435.8237
$ws.Range("K132").Value = 16243.656
$ws.Range("L132").Value = 22307.4711
$ws.Range("M132").Value = -13713.656
$ws.Range("N132").Value = -27367.4711

$ws.Range("H136").Value = 5867.75
$ws.Range("I136").Value = 2104.182
$ws.Range("K136").Value = 6312.545999999999
$ws.Range("M136").Value = -3762.545999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 11127.625
$ws.Range("I18").Value = 29000
$ws.Range("J18").Value = 8574.429
$ws.Range("K18").Value = 29000
$ws.Range("L18").Value = 8574.429
$ws.Range("M18").Value = -28827
$ws.Range("N18").Value = -8920.429

$ws.Range("H41").Value = 15221.75
$ws.Range("I41").Value = 14444
$ws.Range("J41").Value = 15999.5
$ws.Range("K41").Value = 14444
$ws.Range("L41").Value = 15999.5
$ws.Range("M41").Value = -14054
$ws.Range("N41").Value = -16779.5

$ws.Range("H81").Value = 7125.2856
$ws.Range("I81").Value = 5687.5
$ws.Range("J81").Value = 7700.4
$ws.Range("K81").Value = 11375
$ws.Range("L81").Value = 15400.8
$ws.Range("M81").Value = -10314
$ws.Range("N81").Value = -17522.8

$ws.Range("H84").Value = 7125.2856
$ws.Range("I84").Value = 5687.5
$ws.Range("J84").Value = 7700.4
$ws.Range("K84").Value = 56875
$ws.Range("L84").Value = 77004
$ws.Range("M84").Value = -51571
$ws.Range("N84").Value = -87612

$ws.Range("H107").Value = 1086.5625
$ws.Range("I107").Value = 884.6429000000001
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 2653.9287
$ws.Range("L107").Value = 7500
$ws.Range("M107").Value = -733.9287000000004
$ws.Range("N107").Value = -11340

$ws.Range("H108").Value = 74473.5
$ws.Range("J108").Value = 74473.5
$ws.Range("L108").Value = 74473.5
$ws.Range("N108").Value = -82153.5

$ws.Range("H132").Value = 3201.6829
$ws.Range("I132").Value = 2115.2432
$ws.Range("K132").Value = 6345.7296
$ws.Range("M132").Value = -3815.7296
